$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the column headers so the "_old" / "_new" suffixes become the
#    concrete format-version identifiers "_FV2310" / "_FV2404".
# ---------------------------------------------------------------------------
$headers = @{
    "A1" = "Segmentname_FV2310"
    "B1" = "Segmentgruppe_FV2310"
    "C1" = "Segment_FV2310"
    "D1" = "Datenelement_FV2310"
    "E1" = "Segment ID_FV2310"
    "F1" = "Code_FV2310"
    "G1" = "Qualifier_FV2310"
    "H1" = "Beschreibung_FV2310"
    "I1" = "Bedingungsausdruck_FV2310"
    "J1" = "Bedingung_FV2310"
    "L1" = "Segmentname_FV2404"
    "M1" = "Segmentgruppe_FV2404"
    "N1" = "Segment_FV2404"
    "O1" = "Datenelement_FV2404"
    "P1" = "Segment ID_FV2404"
    "Q1" = "Code_FV2404"
    "R1" = "Qualifier_FV2404"
    "S1" = "Beschreibung_FV2404"
    "T1" = "Bedingungsausdruck_FV2404"
    "U1" = "Bedingung_FV2404"
}

foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

# ---------------------------------------------------------------------------
# 2. Freeze the header row (row 1) so it stays visible while scrolling.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.SplitRow = 1
$win.SplitColumn = 0
$win.FreezePanes = $true
$win.SplitRow = 0
$ws.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Turn the used range into a proper Excel Table ("Table1") covering the
#    whole data set (header row + 64 data rows, columns A:U).
# ---------------------------------------------------------------------------
$dataRange = $ws.Range("A1:U65")
$table = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $dataRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$table.Name = "Table1"
